$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "Exportado:" timestamp string (D3)
$ws.Range("D3").Value = "lun. 26/05/2025 21:32"

# 2) Update the last existing data row (row 41, date 07/05/2025) with its new values
$ws.Range("B41").Value = 7
$ws.Range("C41").Value = 121000
$ws.Range("D41").Value = 12

# 3) Grab the table and temporarily hide the totals row so the row it currently
#    occupies (row 42) becomes an ordinary writable cell and the table shrinks
#    back to just the header + existing data (A5:D41). The cells that made up
#    the old totals row keep their old formulas/styles until we overwrite them
#    below.
$lo = $ws.ListObjects.Item(1)
$lo.ShowTotals = $false

# Reset the old totals-row formatting on row 42 back to plain data-row styling
# (date format on column A, general/default on B:D) before we drop new values
# into it, so it matches the rest of the data rows.
$ws.Range("A42").Style = "Normal"
$ws.Range("A42").NumberFormat = "dd/MM/yyyy HH:mm:ss"
$ws.Range("B42:D42").Style = "Normal"

# 4) Append the new daily rows (08/05/2025 .. 26/05/2025) below the existing data
$newData = @(
    @(45785, 10, 155000, 16),
    @(45786, 17, 300100, 29),
    @(45787, 31, 486700, 56),
    @(45788, 14, 223800, 25),
    @(45789, 8, 104342.2, 54),
    @(45790, 67, 376425, 148),
    @(45791, 74, 318930, 152),
    @(45792, 74, 308463, 117),
    @(45793, 72, 536582.5, 151),
    @(45794, 75, 459845, 130),
    @(45795, 42, 357246, 83),
    @(45796, 71, 184170, 105),
    @(45797, 90, 364345.8, 162),
    @(45798, 71, 412825, 134),
    @(45799, 73, 344173, 169),
    @(45800, 83, 413234, 159),
    @(45801, 83, 471428, 144),
    @(45802, 55, 385745, 104),
    @(45803, 62, 162213, 97)
)

$r = 42
foreach ($row in $newData) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# 5) Grow the table back over the new data (A5:D60) and restore the totals row,
#    which lands on row 61 with the same SUBTOTAL formulas.
$lo.Resize($ws.Range("A5:D60"))
$lo.ShowTotals = $true

$ws.Range("A61").NumberFormat = "dd/MM/yyyy HH:mm:ss"
$ws.Range("B61").Formula = "=SUBTOTAL(109,Datos[Cantidad de pedidos])"
$ws.Range("C61").Formula = "=SUBTOTAL(109,Datos[Total $])"
$ws.Range("D61").Formula = "=SUBTOTAL(109,Datos[Cantidad de productos])"

Write-Output "table range:"
Write-Output $lo.Range.Address
